$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New activity-log entries appended after the existing data (rows 2-197).
# Columns: A=User, B=Activity, C=Date, D=Time
$newRows = @(
    @("q", "LOG-IN",  "02/09/24", "11:40:36"),
    @("q", "LOG-IN",  "02/09/24", "11:41:20"),
    @("q", "LOG-IN",  "02/09/24", "11:43:25"),
    @("q", "LOG-OUT", "02/09/24", "11:43:46")
)

$startRow = 198
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    $ws.Range("A$r").Value = $rowData[0]
    $ws.Range("B$r").Value = $rowData[1]

    # Force the date/time-looking strings to be stored as plain text,
    # matching the rest of the sheet (which keeps them as text, not as
    # Excel date/time serial numbers).
    $ws.Range("C$r").NumberFormat = "@"
    $ws.Range("C$r").Value = $rowData[2]

    $ws.Range("D$r").NumberFormat = "@"
    $ws.Range("D$r").Value = $rowData[3]
}
